$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Tipo" column (D) to hold "MAE"
$ws.Range("D1").EntireColumn.Insert()

# Copy style (format) of the existing header cell (e.g. C1) into the new D1 header
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header text
$ws.Range("D1").Value = "MAE"

# Update B/C values and fill in the new D column values
$ws.Range("B2").Value = 0.6829067359667368
$ws.Range("C2").Value = 0.97994124024664
$ws.Range("D2").Value = 0.6227787498567806

$ws.Range("B3").Value = 0.4038961864144212
$ws.Range("C3").Value = 0.9944180453430562
$ws.Range("D3").Value = 0.5245133784119467

$ws.Range("B4").Value = 0.1540378076254267
$ws.Range("C4").Value = 0.9984123651617502
$ws.Range("D4").Value = 0.3359641668615422

$ws.Range("B5").Value = 0.451753926346256
$ws.Range("C5").Value = 0.9973090307564181
$ws.Range("D5").Value = 0.555099691298675

$wb.Save()
